$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChaseData")

# --- Step 1: turn old marker cell G1 ("testcase") into a normal header cell,
#     matching the look of the other header cells (B1:F1), then give it the
#     new header text "AccountID".
$ws.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "AccountID"

# --- Step 2: add the new "SubAccountID" header next to it, same formatting.
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "SubAccountID"

# --- Step 3: re-create the "testcase" marker cell one column further right,
#     reusing A1's look (bold / yellow fill / border) as donor, then strip
#     the left border to match the new style used for I1.
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "testcase"
$ws.Range("I1").Borders.Item(7).LineStyle = -4142   # xlEdgeLeft = xlLineStyleNone

# --- Step 4: add the new data row values under the new headers, matching
#     the plain bordered/centered look used for the rest of row 2 -- but
#     without the centered alignment (per target style xf5).
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "P610"
$ws.Range("G2").HorizontalAlignment = -4108  # clear center -> general, set below precisely
$ws.Range("G2").HorizontalAlignment = 1      # xlHAlignGeneral

$ws.Range("F2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = "Z888"
$ws.Range("H2").HorizontalAlignment = 1      # xlHAlignGeneral

# --- Column widths ---
$ws.Columns.Item(6).ColumnWidth = 12.21875
$ws.Columns.Item(7).ColumnWidth = 10
$ws.Columns.Item(8).ColumnWidth = 13.44140625

# --- Selection ---
$ws.Range("H1").Select()
